$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Scores")

# Update header labels (columns E..J) — shared strings changed
$ws.Range("E1").Value = "GRACE CA"
$ws.Range("F1").Value = "GRACE Exam"
$ws.Range("G1").Value = "AHHBOJ CA"
$ws.Range("H1").Value = "AHHBOJ Exam"
$ws.Range("I1").Value = "JLBO CA"
$ws.Range("J1").Value = "JLBO Exam"

# Update row 2 values
$ws.Range("B2").Value = "HOUR"
$ws.Range("C2").Value = 13
$ws.Range("D2").Value = 59
$ws.Range("E2").Value = 19
$ws.Range("F2").Value = 60
$ws.Range("G2").Value = 11
$ws.Range("H2").Value = 49
$ws.Range("I2").Value = 13
$ws.Range("J2").Value = 14

# Remove row 3 entirely (student data row deleted)
$ws.Rows.Item(3).Delete()
